$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @('U43_01', '43', 'Bối rối', 'Frustrated', 'He is frustrated at this math question.', 'Frustrated at something', 'Adj'),
    @('U43_02', '43', 'Giả định', 'Assumption', 'They make assumptions about me being rude(rằng tôi là người thô lỗ)', 'make an assumption / đưa ra giả định', 'N'),
    @('U43_03', '43', 'Phép thuật', 'magic', 'She made it disappeer like magic', 'like magic / như phép màu', 'N'),
    @('U43_04', '43', 'Cơ hội', 'Chance', 'We met each other by chance / chúng tôi tình cờ gặp nhau', 'By chance / tình cờ', 'N'),
    @('U43_05', '43', 'Yếu tố', 'Factor', 'Have you considered all the important factors', 'an important factor / một yếu tố quan trọng', 'N'),
    @('U43_06', '43', 'Bên trong', 'Inside', 'Your shirt is inside out', 'inside out / lộn ngược phần trong ra ngoài', 'Adj'),
    @('U43_07', '43', 'Sự lo ngại', 'Concern', 'I have lots of concerns about my son', 'concern about something or somebody', 'N'),
    @('U43_08', '43', 'Đánh giá', 'Rate', 'Please rate 5 star for that book on App Store', 'Rate something', 'V'),
    @('U43_09', '43', 'Hoàn toàn', 'Completely', 'He was a completely different person', 'completely different / hoàn toàn khác', 'Adv'),
    @('U43_10', '43', 'Sự tồn tại', 'Existence', 'Dinosaurs are no longer in existence', 'in existence / tồn tại', 'N'),
    @('U43_11', '43', 'Nghiên cứu, xem xét', 'Look into', 'The police will look into the files', 'Look into something', 'V'),
    @('U43_12', '43', 'Kiểm soát', 'Control', 'My father always tries to control me', 'control somebody or something', 'V'),
    @('U43_13', '43', 'Bài viết', 'Article', 'Have you read the news article today?', 'a news article / bài báo', 'N'),
    @('U43_14', '43', 'Tác giả', 'Author', 'Best-selling authors are skilled in writing', 'a best-selling author / tác giả bán chạy nhất', 'N'),
    @('U43_15', '43', 'Đoạn', 'Passage', 'please read the message on page 24.', 'read a passage / đọc một đoạn', 'N'),
    @('U43_16', '43', 'Thuyết phục', 'Convince', 'Please convince her to help me', 'convince somebody', 'V'),
    @('U43_17', '43', 'Phần', 'Part', 'I am a part of the team', 'part of something', 'N'),
    @('U43_18', '43', 'Vừa đủ', 'Adequate', 'She show an adequate amout of skill.', 'An adequate amount / một lượng vừa đủ', 'Adj'),
    @('U43_19', '43', 'Chọn', 'Pick', 'The team pucked him to be the leader', 'pick somebody or something', 'V'),
    @('U43_20', '43', 'Thành thật, trung thực', 'Honest', 'Please be honest about your past.', 'honest about something / thành thật về một điều gì đó', 'Adj'),
    @('U43_21', '43', 'Thất vọng', 'Disappointed', 'I am disappointed with what you did', 'disappointed with', 'Adj'),
    @('U43_22', '43', 'Kết cấu', 'Structure', 'The structure of the company is complex', 'a structure of something', 'N'),
    @('U43_23', '43', 'Tập phim', 'Episode', 'They broke up in the last episode', 'in an episode / trong tập phim', 'N'),
    @('U43_24', '43', 'Thay thế', 'Replace', 'The school replaced the old building', 'replace somebody or something', 'V'),
    @('U43_25', '43', 'Nguyên bản, độc đáo', 'Original', 'He came up with an original idea.', 'an original idea / ý tưởng độc đáo', 'Adj'),
    @('U43_26', '43', 'Sự phân tích', 'Analysis', 'We need an expert in data analysis', 'data analysis / việc phân tích dữ liệu', 'N'),
    @('U43_27', '43', 'Công lý', 'Justice', 'Justice will be served when he receives punishment(hình phạt)', 'justice is served / công lý được thực thi', 'N'),
    @('U43_28', '43', 'Đại diện cho, ủng hộ cho', 'Stand for', 'We stand for equality', 'stand for something / ủng hộ cho một cái gì đó', 'V'),
    @('U43_29', '43', 'Kết luận', 'Conclusion', 'In conclusion, this book is good', 'in conclusion / kết luận là', 'N'),
    @('U43_30', '43', 'Khắp, suốt', 'Throughout', 'The school opens throughout the year', 'throughout the year / suốt cả năm', 'Adv')
)

$startRow = 1262
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = [int]$row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

$ws.Range("C1292").Select()
